$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestData")

# Update the username/email and name values on the TestData sheet
$ws1.Range("B2").Value = "himabejo1@gmail.com"
$ws1.Range("D2").Value = "Hima"

# Update the active selection on the TestData sheet
$ws1.Range("E2").Select() | Out-Null
